$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2073863636363636
$ws.Range("C2").Value = 0.5397727272727273
$ws.Range("J2").Value = 0.01420454545454545
$ws.Range("P2").Value = 0.1505681818181818
$ws.Range("S2").Value = 0.08806818181818182

# Row 3
$ws.Range("B3").Value = 0.005
$ws.Range("C3").Value = 0.05
$ws.Range("J3").Value = 0.025
$ws.Range("P3").Value = 0.745
$ws.Range("S3").Value = 0.175

# Row 4
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.7727272727272727
$ws.Range("S4").Value = 0.2045454545454546

# Row 6
$ws.Range("B6").Value = 0.0931174089068826
$ws.Range("D6").Value = 0.008097165991902834
$ws.Range("F6").Value = 0.04453441295546558
$ws.Range("J6").Value = 0.2591093117408907
$ws.Range("O6").Value = 0.03238866396761134
$ws.Range("Q6").Value = 0.1659919028340081
$ws.Range("R6").Value = 0.0931174089068826
$ws.Range("S6").Value = 0.3036437246963563

# Row 7
$ws.Range("B7").Value = 0.08771929824561403
$ws.Range("D7").Value = 0.01754385964912281
$ws.Range("F7").Value = 0.07602339181286549
$ws.Range("J7").Value = 0.1052631578947368
$ws.Range("O7").Value = 0.01169590643274854
$ws.Range("Q7").Value = 0.2046783625730994
$ws.Range("R7").Value = 0.06432748538011696
$ws.Range("S7").Value = 0.4327485380116959

# Row 8
$ws.Range("B8").Value = 0.1185344827586207
$ws.Range("D8").Value = 0.02801724137931035
$ws.Range("E8").Value = 0.004310344827586207
$ws.Range("F8").Value = 0.07758620689655173
$ws.Range("J8").Value = 0.1120689655172414
$ws.Range("O8").Value = 0.03017241379310345
$ws.Range("Q8").Value = 0.2133620689655172
$ws.Range("R8").Value = 0.08836206896551724
$ws.Range("S8").Value = 0.3275862068965517

# Row 9
$ws.Range("B9").Value = 0.1317365269461078
$ws.Range("D9").Value = 0.005988023952095809
$ws.Range("F9").Value = 0.08383233532934131
$ws.Range("J9").Value = 0.155688622754491
$ws.Range("O9").Value = 0.005988023952095809
$ws.Range("Q9").Value = 0.1736526946107785
$ws.Range("R9").Value = 0.125748502994012
$ws.Range("S9").Value = 0.3173652694610778

# Row 10
$ws.Range("B10").Value = 0.1262214983713355
$ws.Range("D10").Value = 0.02361563517915309
$ws.Range("E10").Value = 0.003257328990228013
$ws.Range("F10").Value = 0.07491856677524431
$ws.Range("J10").Value = 0.1311074918566775
$ws.Range("O10").Value = 0.02117263843648208
$ws.Range("Q10").Value = 0.1978827361563518
$ws.Range("R10").Value = 0.1034201954397394
$ws.Range("S10").Value = 0.3184039087947882

# Row 11
$ws.Range("G11").Value = 0.1330798479087452
$ws.Range("J11").Value = 0.08745247148288973
$ws.Range("K11").Value = 0.2053231939163498
$ws.Range("L11").Value = 0.5627376425855514
$ws.Range("S11").Value = 0.01140684410646388

# Row 12
$ws.Range("G12").Value = 0.7516339869281046
$ws.Range("J12").Value = 0.1764705882352941
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.0392156862745098
$ws.Range("S12").Value = 0.0261437908496732

# Row 13
$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.2564102564102564
$ws.Range("S13").Value = 0.05128205128205128

# Row 15
$ws.Range("F15").Value = 0.02459016393442623
$ws.Range("H15").Value = 0.1516393442622951
$ws.Range("I15").Value = 0.03688524590163934
$ws.Range("J15").Value = 0.3442622950819672
$ws.Range("K15").Value = 0.06967213114754098
$ws.Range("M15").Value = 0.01229508196721311
$ws.Range("O15").Value = 0.09426229508196721
$ws.Range("S15").Value = 0.2663934426229508

# Row 16
$ws.Range("F16").Value = 0.008695652173913044
$ws.Range("H16").Value = 0.208695652173913
$ws.Range("I16").Value = 0.06086956521739131
$ws.Range("J16").Value = 0.3695652173913043
$ws.Range("K16").Value = 0.08695652173913043
$ws.Range("M16").Value = 0.02173913043478261
$ws.Range("N16").Value = 0.004347826086956522
$ws.Range("O16").Value = 0.07391304347826087
$ws.Range("S16").Value = 0.1652173913043478

# Row 17
$ws.Range("F17").Value = 0.02237136465324385
$ws.Range("H17").Value = 0.1968680089485459
$ws.Range("I17").Value = 0.08724832214765101
$ws.Range("J17").Value = 0.4205816554809844
$ws.Range("K17").Value = 0.1006711409395973
$ws.Range("M17").Value = 0.01565995525727069
$ws.Range("N17").Value = 0.002237136465324385
$ws.Range("O17").Value = 0.06935123042505593
$ws.Range("S17").Value = 0.08501118568232663

# Row 18
$ws.Range("F18").Value = 0.02690582959641256
$ws.Range("H18").Value = 0.2197309417040359
$ws.Range("I18").Value = 0.1121076233183857
$ws.Range("J18").Value = 0.336322869955157
$ws.Range("K18").Value = 0.09865470852017937
$ws.Range("M18").Value = 0.008968609865470852
$ws.Range("O18").Value = 0.08520179372197309
$ws.Range("S18").Value = 0.1121076233183857

# Row 19
$ws.Range("F19").Value = 0.01931922723091076
$ws.Range("H19").Value = 0.2207911683532659
$ws.Range("I19").Value = 0.07359705611775529
$ws.Range("J19").Value = 0.3946642134314627
$ws.Range("K19").Value = 0.09751609935602576
$ws.Range("M19").Value = 0.0202391904323827
$ws.Range("N19").Value = 0.002759889604415824
$ws.Range("O19").Value = 0.07543698252069918
$ws.Range("S19").Value = 0.09567617295308188

